{"js": "const replacements = [\n  [\"631\u00d79=\", \"871\u00d78=\"],\n  [\"297\u00d72=\", \"576\u00d76=\"],\n  [\"454\u00d73=\", \"296\u00d79=\"],\n  [\"928\u00d77=\", \"925\u00d76=\"],\n  [\"486\u00d75=\", \"863\u00d75=\"],\n  [\"564\u00d79=\", \"348\u00d72=\"],\n  [\"645\u00d79=\", \"442\u00d76=\"],\n  [\"702\u00d76=\", \"242\u00d75=\"],\n  [\"361\u00d77=\", \"592\u00d77=\"],\n  [\"555\u00d74=\", \"990\u00d76=\"],\n  [\"256\u00d73=\", \"277\u00d72=\"],\n  [\"504\u00d72=\", \"438\u00d72=\"],\n  [\"401\u00d74=\", \"781\u00d77=\"],\n  [\"413\u00d75=\", \"838\u00d75=\"],\n  [\"903\u00d77=\", \"370\u00d78=\"],\n  [\"979\u00d74=\", \"405\u00d79=\"],\n  [\"657\u00d72=\", \"454\u00d76=\"],\n  [\"387\u00d74=\", \"417\u00d79=\"],\n  [\"885\u00d75=\", \"847\u00d73=\"],\n  [\"972\u00d78=\", \"633\u00d74=\"],\n  [\"986\u00d73=\", \"300\u00d75=\"],\n  [\"870\u00d79=\", \"779\u00d75=\"],\n  [\"193\u00d72=\", \"874\u00d78=\"],\n  [\"153\u00d75=\", \"152\u00d75=\"],\n  [\"645\u00d76=\", \"849\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"631\u00d79=\", \"871\u00d78=\")\n  ,@(\"297\u00d72=\", \"576\u00d76=\")\n  ,@(\"454\u00d73=\", \"296\u00d79=\")\n  ,@(\"928\u00d77=\", \"925\u00d76=\")\n  ,@(\"486\u00d75=\", \"863\u00d75=\")\n  ,@(\"564\u00d79=\", \"348\u00d72=\")\n  ,@(\"645\u00d79=\", \"442\u00d76=\")\n  ,@(\"702\u00d76=\", \"242\u00d75=\")\n  ,@(\"361\u00d77=\", \"592\u00d77=\")\n  ,@(\"555\u00d74=\", \"990\u00d76=\")\n  ,@(\"256\u00d73=\", \"277\u00d72=\")\n  ,@(\"504\u00d72=\", \"438\u00d72=\")\n  ,@(\"401\u00d74=\", \"781\u00d77=\")\n  ,@(\"413\u00d75=\", \"838\u00d75=\")\n  ,@(\"903\u00d77=\", \"370\u00d78=\")\n  ,@(\"979\u00d74=\", \"405\u00d79=\")\n  ,@(\"657\u00d72=\", \"454\u00d76=\")\n  ,@(\"387\u00d74=\", \"417\u00d79=\")\n  ,@(\"885\u00d75=\", \"847\u00d73=\")\n  ,@(\"972\u00d78=\", \"633\u00d74=\")\n  ,@(\"986\u00d73=\", \"300\u00d75=\")\n  ,@(\"870\u00d79=\", \"779\u00d75=\")\n  ,@(\"193\u00d72=\", \"874\u00d78=\")\n  ,@(\"153\u00d75=\", \"152\u00d75=\")\n  ,@(\"645\u00d76=\", \"849\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    Write-Output \"WARNING: replacement not found for $oldText\"\n  }\n}"}
